$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customerlogin")

$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"
